$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Translate Publisher and Contact values from German to English
$wsMeta.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$wsMeta.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Add a Description value next to the "Description" property (row 12)
$wsMeta.Range("B12").Value = "Extensible description of the application perspective of a ConsentPolicy (e.g. selected data view 'IDAT') "
